$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (masthead) updates ---
$ws.Range("A8").Value = "Volume 31   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/1/2024  Through  4/7/2024"

# Donor cells holding the canonical "0" / "***.*" text+style (row 14 - Murder - is untouched by this edit)
$zeroDonor = $ws.Range("C14")
$naDonor = $ws.Range("E14")

# --- Row 15 ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("N15").Value = -12.5

# --- Row 16 ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -5
$ws.Range("L16").Value = -15.555555555555
$ws.Range("M16").Value = -20.833333333333
$ws.Range("N16").Value = -88.622754491018

# --- Row 17 ---
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 4.347826086956
$ws.Range("I17").Value = 81
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 42.105263157894
$ws.Range("L17").Value = 3.846153846153
$ws.Range("M17").Value = 84.090909090909
$ws.Range("N17").Value = -36.71875

# --- Row 18 ---
$ws.Range("C18").Value = 3
$zeroDonor.Copy($ws.Range("D18"))
$naDonor.Copy($ws.Range("E18"))
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 14.285714285714
$ws.Range("L18").Value = -37.662337662337
$ws.Range("M18").Value = -34.246575342465
$ws.Range("N18").Value = -83.561643835616

# --- Row 19 ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -46.153846153846
$ws.Range("I19").Value = 140
$ws.Range("J19").Value = 154
$ws.Range("K19").Value = -9.090909090909
$ws.Range("L19").Value = -49.27536231884
$ws.Range("M19").Value = 72.839506172839
$ws.Range("N19").Value = 37.254901960784

# --- Row 20 ---
$ws.Range("C20").Value = 6
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 6
$ws.Range("E20").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -13.333333333333
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 29.411764705882
$ws.Range("M20").Value = 12.820512820512
$ws.Range("N20").Value = -85.185185185185

# --- Row 21 ---
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 29.629629629629
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -19.298245614035
$ws.Range("I21").Value = 358
$ws.Range("J21").Value = 345
$ws.Range("K21").Value = 3.768115942028
$ws.Range("L21").Value = -30.754352030947
$ws.Range("M21").Value = 22.602739726027
$ws.Range("N21").Value = -69.270386266094

# --- Row 22 ---
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666

# --- Row 24 ---
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -10.526315789473
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 10.843373493975
$ws.Range("I24").Value = 268
$ws.Range("J24").Value = 279
$ws.Range("K24").Value = -3.942652329749
$ws.Range("L24").Value = -22.766570605187
$ws.Range("M24").Value = 37.435897435897

# --- Row 25 ---
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -14.545454545454
$ws.Range("I25").Value = 130
$ws.Range("J25").Value = 166
$ws.Range("K25").Value = -21.686746987951
$ws.Range("L25").Value = -33.333333333333

# --- Row 26 ---
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 80
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 6.666666666666
$ws.Range("I26").Value = 118
$ws.Range("J26").Value = 113
$ws.Range("K26").Value = 4.424778761061
$ws.Range("L26").Value = 9.259259259259
$ws.Range("M26").Value = -30.588235294117

# --- Row 27 ---
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$zeroDonor.Copy($ws.Range("D27"))
$naDonor.Copy($ws.Range("E27"))
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50

# --- Row 28 ---
$zeroDonor.Copy($ws.Range("C28"))
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -6.666666666666

# --- Row 29 ---
$ws.Range("G29").Value = 1

# --- Row 30 ---
$ws.Range("G30").Value = 1
